$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2 and 3 in column B
$ws.Range("B2").Value = 8
$ws.Range("B3").Value = 7

# Delete row 4 entirely (shifts cells up, removing the row from the sheet)
$ws.Rows.Item(4).Delete()
